$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.788.14"
$ws.Range("E2").Value = "  +4.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.52"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.93"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  -0.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  +3.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4039"
$ws.Range("E8").Value = "  +5.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.72"
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08081"
$ws.Range("E10").Value = "  +2.34%  "
$ws.Range("E11").Value = "  +4.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.29"
$ws.Range("E12").Value = "  +5.87%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.884.93"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.071"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.322"
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.86"
$ws.Range("E16").Value = "  +1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001050"
$ws.Range("E18").Value = "  +2.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06620"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9989"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.795.25"
$ws.Range("E22").Value = "  +4.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.525"
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.264"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.099.72"
$ws.Range("E26").Value = "  +1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.70"
$ws.Range("E27").Value = "  +3.12%  "
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.144"
$ws.Range("E29").Value = "  +4.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.518"
$ws.Range("E30").Value = "  +4.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.31"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9949"
$ws.Range("E32").Value = "  +5.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09566"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.659"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  +5.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.393"
$ws.Range("E36").Value = "  +2.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06200"
$ws.Range("E37").Value = "  +4.57%  "
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.526"
$ws.Range("E39").Value = "  +6.16%  "
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5988"
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1896"
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.40"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.269"
$ws.Range("E45").Value = "  +0.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5608"
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.26"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.971"
$ws.Range("E48").Value = "  +5.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07223"
$ws.Range("E49").Value = "  +9.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.119"
$ws.Range("E50").Value = "  +14.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.71"
$ws.Range("E51").Value = "  +1.90%  "
